$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the paragraph that contains the existing Streamlit hyperlink
# ("https://p03---chla-pred-hl.streamlit.app/"). The new GitHub-link
# paragraph must be inserted immediately after it.
# ---------------------------------------------------------------------------
$streamlitUrl = "https://p03---chla-pred-hl.streamlit.app/"
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd() -eq $streamlitUrl) {
        $anchorPara = $cand
        break
    }
}
if ($null -eq $anchorPara) {
    throw "Could not locate the paragraph containing the Streamlit hyperlink."
}

# ---------------------------------------------------------------------------
# Insert a brand-new, empty paragraph right after it.
# ---------------------------------------------------------------------------
$anchorRange = $anchorPara.Range
$anchorRange.Collapse(0)
$anchorRange.InsertParagraphAfter()

$newPara = $anchorPara.Next()
$newRange = $newPara.Range

# ---------------------------------------------------------------------------
# Fill the new paragraph with a hyperlink run pointing at the GitHub repo,
# via a WordOpenXML package fragment (adds both the run/hyperlink markup and
# the External hyperlink relationship in one shot).
# ---------------------------------------------------------------------------
$githubUrl = "https://github.com/Hlissner31/P03---CHLA-PRED"

$packageXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
<Relationship Id="rIdGhLink" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="__GITHUB_URL__" TargetMode="External"/>
</Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p>
<w:hyperlink r:id="rIdGhLink" w:history="1">
<w:r>
<w:t>__GITHUB_URL__</w:t>
</w:r>
</w:hyperlink>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$packageXml = $packageXml.Replace("__GITHUB_URL__", $githubUrl)
$newRange.InsertXML($packageXml) | Out-Null

# ---------------------------------------------------------------------------
# InsertXML drops w:rStyle references, so re-apply the "Hyperlink" character
# style to the inserted run explicitly, matching the existing hyperlink run.
# ---------------------------------------------------------------------------
$insertedPara = $anchorPara.Next()
$runRange = $insertedPara.Range
$runRange.MoveEnd(1, -1) | Out-Null
$runRange.Style = $d.Styles.Item("Hyperlink")

Write-Host "Inserted paragraph text:" $insertedPara.Range.Text
Write-Host "Total hyperlinks in document:" $d.Hyperlinks.Count
